# Module 1 Quiz - CAGR worksheet: add a second block (rows 13-22) that
# restates each year's return as "1 + return" (growth factor), and a
# summary row (23) that backs out the CAGR of each column via GEOMEAN.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Growth-factor block: rows 13-22 mirror rows 2-11 -----------------
# Column A: just carries the year forward from the row above (2 through 11)
$ws.Range("A13").Formula = "=A2"
$ws.Range("A14:A21").Formula = "=A3"
$ws.Range("A22").Formula = "=A11"

# Columns B:E: 1 + the corresponding annual return, one row at a time so
# each row becomes its own shared-formula group (matches a manual fill-down
# per row rather than one single fill across the whole block).
$ws.Range("B13:E13").Formula = "=B2+1"
$ws.Range("B14:E14").Formula = "=B3+1"
$ws.Range("B15:E15").Formula = "=B4+1"
$ws.Range("B16:E16").Formula = "=B5+1"
$ws.Range("B17:E17").Formula = "=B6+1"
$ws.Range("B18:E18").Formula = "=B7+1"
$ws.Range("B19:E19").Formula = "=B8+1"
$ws.Range("B20:E20").Formula = "=B9+1"
$ws.Range("B21:E21").Formula = "=B10+1"
$ws.Range("B22:E22").Formula = "=B11+1"

$ws.Range("B13:E22").NumberFormat = "0.000%"

# --- Summary row: geometric mean of the growth factors, minus 1 -------
$ws.Range("B23:E23").Formula = "=GEOMEAN(B13:B22)-1"
$ws.Range("B23:E23").NumberFormat = '_(* #,##0.0000_);_(* \(#,##0.0000\);_(* "-"??_);_(@_)'

# --- Cosmetics ----------------------------------------------------------
# Column B ("IBM") grew a bit wider once it picked up the longer/boldened
# GEOMEAN summary figures underneath it - approximate Excel's best-fit here.
$ws.Columns("B").ColumnWidth = 21.5

$ws.Range("H12").Select() | Out-Null
